$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.343.09'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.857.37'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.83%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.00'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4609'
$ws.Range('E7').Value = '  -0.94%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3706'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07325'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8814'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.87'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07799'
$ws.Range('E12').Value = '  -1.48%  '
$ws.Range('D13').Value = '1.830.13'
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.391'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.545'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.84'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000009014'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').Value = '27.358.58'
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('D24').Value = '2.142.35'
$ws.Range('E24').Value = '  +4.35%  '
$ws.Range('E25').Value = '  +4.85%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '151.98'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.36'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.074'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.115'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '116.18'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08826'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7704'
$ws.Range('E32').Value = '  +5.64%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.018'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('E34').Value = '  +3.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.491'
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.636'
$ws.Range('E36').Value = '  +5.92%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01961'
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.078'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05228'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.950'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.029'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5145'
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1640'
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.404'
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.4839'
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.30'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.9999'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '103.36'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.654'
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06216'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '65.93'
$ws.Range('E51').Value = '  +2.11%  '
